# Apply the "Branch wise aging stock status table" update.
# The underlying edit reorders several item rows within a few brand
# groups (Dinafex, Etorix, Flucloxin, Ketonic, Zithrox) and corrects a
# handful of numeric stock figures (Geminox MTD Sales Target, and the
# Sk-Mox aging-stock row). Re-create the resulting cell values directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Dinafex group (rows 3-5): item names rotate, TP (BB) follows them ---
$ws.Range("C3").Value = "Dinafex 120mg Tablet"
$ws.Range("BB3").Value = 179.91

$ws.Range("C4").Value = "Dinafex 60mg Tablet"
$ws.Range("BB4").Value = 78.70999999999999

$ws.Range("C5").Value = "Dinafex 180mg Tablet"
$ws.Range("BB5").Value = 224.89

# --- Etorix group (rows 7-9): item names + UOM rotate, TP follows ---
$ws.Range("C7").Value = "Etorix 60mg Tablet - 40's"
$ws.Range("D7").Value = "40's"

$ws.Range("C8").Value = "Etorix 90mg Tablet"
$ws.Range("D8").Value = "30's"
$ws.Range("BB8").Value = 269.87

$ws.Range("C9").Value = "Etorix 120mg Tablet"
$ws.Range("D9").Value = "20's"
$ws.Range("BB9").Value = 209.9

# --- Flucloxin group (rows 11-12): item names + UOM swap, TP follows ---
$ws.Range("C11").Value = "Flucloxin 500mg Capsule - 36's"
$ws.Range("D11").Value = "36 's"
$ws.Range("BB11").Value = 284.21

$ws.Range("C12").Value = "Flucloxin 500mg Capsule"
$ws.Range("D12").Value = "30 's"
$ws.Range("BB12").Value = 237.74

# --- Geminox (row 13): MTD Sales Target correction ---
$ws.Range("H13").Value = 103

# --- Ketonic group (rows 14 & 16): item names + UOM swap, TP follows ---
$ws.Range("C14").Value = "Ketonic 30mg IM/IV Injection - 4's"
$ws.Range("D14").Value = "4's"
$ws.Range("BB14").Value = 165.41

$ws.Range("C16").Value = "Ketonic 10mg Tablet"
$ws.Range("D16").Value = "20's"
$ws.Range("BB16").Value = 150.38

# --- Sk-Mox (row 23): aging stock status figures corrected ---
$ws.Range("E23").Value = 0
$ws.Range("L23").Value = 25
$ws.Range("N23").Value = 167
$ws.Range("O23").Value = 177
$ws.Range("T23").Value = 177
$ws.Range("AL23").Value = 21
$ws.Range("BC23").Value = 0
$ws.Range("BD23").Value = 0

# --- Zithrox group (rows 25-27): item names + UOM rotate, TP follows ---
$ws.Range("C25").Value = "Zithrox 15ml Suspension"
$ws.Range("D25").Value = "15 ml"
$ws.Range("BB25").Value = 71.95999999999999

$ws.Range("C26").Value = "Zithrox 30ml Dry Suspension"
$ws.Range("D26").Value = "30ml"
$ws.Range("BB26").Value = 97.45

$ws.Range("C27").Value = "Zithrox 500mg Tablet"
$ws.Range("D27").Value = "6 's"
$ws.Range("BB27").Value = 136.83
